# Applies the "Update file from SharePoint" edit: appends 26 new rows of
# mailbox-email metadata (received_date_time, subject, sender, attachment)
# pulled in since the prior save, then grows Table1 / the sheet dimension
# to cover them, and nudges column D's width to match the refreshed
# Excel Online autosize pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2038, 1).Value = "2025-07-24T19:16:15+00:00"
$ws.Cells.Item(2038, 2).Value = "EXTERNAL:- BBC Radio nan Gàidheal - Wk31 - 2025-08-03 - Sunday"
$ws.Cells.Item(2038, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2038, 4).Value = $true

$ws.Cells.Item(2039, 1).Value = "2025-07-24T19:16:15+00:00"
$ws.Cells.Item(2039, 2).Value = "EXTERNAL:- BBC Radio 1Xtra - Wk31 - 2025-08-03 - Sunday"
$ws.Cells.Item(2039, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2039, 4).Value = $true

$ws.Cells.Item(2040, 1).Value = "2025-07-24T19:16:13+00:00"
$ws.Cells.Item(2040, 2).Value = "EXTERNAL:- BBC Radio 1 - Wk30 - 2025-07-28 - Monday"
$ws.Cells.Item(2040, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2040, 4).Value = $true

$ws.Cells.Item(2041, 1).Value = "2025-07-24T19:16:07+00:00"
$ws.Cells.Item(2041, 2).Value = "EXTERNAL:- BBC Radio 1 - Wk30 - 2025-07-30 - Wednesday"
$ws.Cells.Item(2041, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2041, 4).Value = $true

$ws.Cells.Item(2042, 1).Value = "2025-07-24T19:15:59+00:00"
$ws.Cells.Item(2042, 2).Value = "EXTERNAL:- BBC Radio 4 Extra - Wk30 - 2025-07-29 - Tuesday"
$ws.Cells.Item(2042, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2042, 4).Value = $true

$ws.Cells.Item(2043, 1).Value = "2025-07-24T19:15:58+00:00"
$ws.Cells.Item(2043, 2).Value = "EXTERNAL:- BBC Radio 4 FM - Wk30 - 2025-08-01 - Friday"
$ws.Cells.Item(2043, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2043, 4).Value = $true

$ws.Cells.Item(2044, 1).Value = "2025-07-24T19:15:58+00:00"
$ws.Cells.Item(2044, 2).Value = "EXTERNAL:- BBC Two HD - Wk31 - 2025-08-03 - Sunday"
$ws.Cells.Item(2044, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2044, 4).Value = $true

$ws.Cells.Item(2045, 1).Value = "2025-07-24T19:15:56+00:00"
$ws.Cells.Item(2045, 2).Value = "EXTERNAL:- BBC Radio 1 - Wk30 - 2025-07-31 - Thursday"
$ws.Cells.Item(2045, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2045, 4).Value = $true

$ws.Cells.Item(2046, 1).Value = "2025-07-24T19:15:56+00:00"
$ws.Cells.Item(2046, 2).Value = "EXTERNAL:- BBC Two HD - Wk31 - 2025-08-03 - Sunday"
$ws.Cells.Item(2046, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2046, 4).Value = $true

$ws.Cells.Item(2047, 1).Value = "2025-07-24T19:15:47+00:00"
$ws.Cells.Item(2047, 2).Value = "EXTERNAL:- BBC Radio 1 - Wk30 - 2025-07-29 - Tuesday"
$ws.Cells.Item(2047, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2047, 4).Value = $true

$ws.Cells.Item(2048, 1).Value = "2025-07-24T20:15:50+00:00"
$ws.Cells.Item(2048, 2).Value = "EXTERNAL:- BBC Radio 4 FM - Wk29 - 2025-07-24 - Thursday"
$ws.Cells.Item(2048, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2048, 4).Value = $true

$ws.Cells.Item(2049, 1).Value = "2025-07-24T20:01:30+00:00"
$ws.Cells.Item(2049, 2).Value = "EXTERNAL:- BBC World Service UK Schedule - Wk34 - 2025-08-28 - Thursday"
$ws.Cells.Item(2049, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2049, 4).Value = $true

$ws.Cells.Item(2050, 1).Value = "2025-07-24T20:01:15+00:00"
$ws.Cells.Item(2050, 2).Value = "EXTERNAL:- BBC Asian Network - Wk33 - 2025-08-22 - Friday"
$ws.Cells.Item(2050, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2050, 4).Value = $true

$ws.Cells.Item(2051, 1).Value = "2025-07-24T20:01:11+00:00"
$ws.Cells.Item(2051, 2).Value = "EXTERNAL:- BBC Radio 3 - Wk33 - 2025-08-19 - Tuesday"
$ws.Cells.Item(2051, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2051, 4).Value = $true

$ws.Cells.Item(2052, 1).Value = "2025-07-24T20:01:08+00:00"
$ws.Cells.Item(2052, 2).Value = "EXTERNAL:- BBC Radio 3 - Wk33 - 2025-08-21 - Thursday"
$ws.Cells.Item(2052, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2052, 4).Value = $true

$ws.Cells.Item(2053, 1).Value = "2025-07-24T20:01:07+00:00"
$ws.Cells.Item(2053, 2).Value = "EXTERNAL:- BBC Radio 1Xtra - Wk33 - 2025-08-20 - Wednesday"
$ws.Cells.Item(2053, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2053, 4).Value = $true

$ws.Cells.Item(2054, 1).Value = "2025-07-24T20:01:03+00:00"
$ws.Cells.Item(2054, 2).Value = "EXTERNAL:- BBC Asian Network - Wk33 - 2025-08-20 - Wednesday"
$ws.Cells.Item(2054, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2054, 4).Value = $true

$ws.Cells.Item(2055, 1).Value = "2025-07-24T20:01:03+00:00"
$ws.Cells.Item(2055, 2).Value = "EXTERNAL:- BBC Radio 6 Music - Wk33 - 2025-08-19 - Tuesday"
$ws.Cells.Item(2055, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2055, 4).Value = $true

$ws.Cells.Item(2056, 1).Value = "2025-07-24T20:01:01+00:00"
$ws.Cells.Item(2056, 2).Value = "EXTERNAL:- BBC Radio 2 - Wk33 - 2025-08-19 - Tuesday"
$ws.Cells.Item(2056, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2056, 4).Value = $true

$ws.Cells.Item(2057, 1).Value = "2025-07-24T20:00:59+00:00"
$ws.Cells.Item(2057, 2).Value = "EXTERNAL:- BBC Asian Network - Wk33 - 2025-08-18 - Monday"
$ws.Cells.Item(2057, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2057, 4).Value = $true

$ws.Cells.Item(2058, 1).Value = "2025-07-24T20:00:54+00:00"
$ws.Cells.Item(2058, 2).Value = "EXTERNAL:- BBC Radio 2 - Wk33 - 2025-08-17 - Sunday"
$ws.Cells.Item(2058, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2058, 4).Value = $true

$ws.Cells.Item(2059, 1).Value = "2025-07-24T20:00:52+00:00"
$ws.Cells.Item(2059, 2).Value = "EXTERNAL:- BBC Radio 4 FM - Wk32 - 2025-08-13 - Wednesday"
$ws.Cells.Item(2059, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2059, 4).Value = $true

$ws.Cells.Item(2060, 1).Value = "2025-07-24T20:00:51+00:00"
$ws.Cells.Item(2060, 2).Value = "EXTERNAL:- BBC ALBA - Wk31 - 2025-08-05 - Tuesday"
$ws.Cells.Item(2060, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2060, 4).Value = $true

$ws.Cells.Item(2061, 1).Value = "2025-07-24T20:00:43+00:00"
$ws.Cells.Item(2061, 2).Value = "EXTERNAL:- BBC Radio 5 Sports Extra - Wk32 - 2025-08-12 - Tuesday"
$ws.Cells.Item(2061, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2061, 4).Value = $true

$ws.Cells.Item(2062, 1).Value = "2025-07-24T20:00:38+00:00"
$ws.Cells.Item(2062, 2).Value = "EXTERNAL:- BBC Radio 1Xtra - Wk32 - 2025-08-10 - Sunday"
$ws.Cells.Item(2062, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2062, 4).Value = $true

$ws.Cells.Item(2063, 1).Value = "2025-07-24T20:00:36+00:00"
$ws.Cells.Item(2063, 2).Value = "EXTERNAL:- BBC Radio 1Xtra - Wk31 - 2025-08-05 - Tuesday"
$ws.Cells.Item(2063, 3).Value = "pressportal@bbc.co.uk"
$ws.Cells.Item(2063, 4).Value = $true

# Grow the table (ListObject) to cover the newly populated rows so its
# `ref`/`autoFilter` stay in sync with the sheet data, matching Table1's
# new A1:D2063 extent.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:D" + $lastRow))

# Column D was re-autofit by Excel Online on this save (12.57 -> 12.43
# characters); 11.67 is the closest COM-addressable width that lands on
# the same 1/6-character XML grid as the target.
$ws.Columns.Item(4).ColumnWidth = 11.67
